$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Cards/kartapolaka.jpg"
